# Add: Calculating Sum and Adding Data
#
# Adds a new row 4 under the GRADES table with the average of each
# sub-column (sum of the two data rows divided by 2), expressed as
# formulas so Excel keeps them live.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Formula = "=SUM(B2:B3)/2"
$ws.Range("C4").Formula = "=SUM(C2:C3)/2"
$ws.Range("D4").Formula = "=SUM(D2:D3)/2"
$ws.Range("E4").Formula = "=SUM(E2:E3)/2"
